# Update "想去人数" (F column) figures across all four sheets, and
# mark one ticket-price cell (G13 on 演出) as unavailable ("不可售").

$wb = $excel.ActiveWorkbook

# 展览 (sheet 1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 8060
$ws1.Range("F6").Value = 4934
$ws1.Range("F7").Value = 7252
$ws1.Range("F8").Value = 886
$ws1.Range("F9").Value = 216
$ws1.Range("F10").Value = 84
$ws1.Range("F11").Value = 1204
$ws1.Range("F13").Value = 212
$ws1.Range("F14").Value = 583
$ws1.Range("F16").Value = 49
$ws1.Range("F17").Value = 254
$ws1.Range("F20").Value = 1309
$ws1.Range("F21").Value = 1281
$ws1.Range("F23").Value = 49
$ws1.Range("F24").Value = 1299
$ws1.Range("F30").Value = 235
$ws1.Range("F33").Value = 20
$ws1.Range("F35").Value = 151
$ws1.Range("F38").Value = 585
$ws1.Range("F41").Value = 71
$ws1.Range("F42").Value = 126
$ws1.Range("F43").Value = 450
$ws1.Range("F45").Value = 636
$ws1.Range("F46").Value = 178
$ws1.Range("F47").Value = 35

# 演出 (sheet 2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 34
$ws2.Range("F10").Value = 14
$ws2.Range("F12").Value = 1739
$ws2.Range("G13").Value = "不可售"
$ws2.Range("F20").Value = 152
$ws2.Range("F35").Value = 2
$ws2.Range("F36").Value = 133
$ws2.Range("F44").Value = 87
$ws2.Range("F46").Value = 13

# 本地生活 (sheet 3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F6").Value = 715
$ws3.Range("F7").Value = 230
$ws3.Range("F9").Value = 1847
$ws3.Range("F10").Value = 2748

# 全部类型 (sheet 4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 715
$ws4.Range("F7").Value = 8060
$ws4.Range("F8").Value = 230
$ws4.Range("F9").Value = 4934
$ws4.Range("F10").Value = 7252
$ws4.Range("F11").Value = 14
$ws4.Range("F12").Value = 886
$ws4.Range("F14").Value = 216
$ws4.Range("F15").Value = 1847
$ws4.Range("F16").Value = 2748
$ws4.Range("F18").Value = 84
$ws4.Range("F19").Value = 1204
$ws4.Range("F20").Value = 212
$ws4.Range("F21").Value = 583
$ws4.Range("F22").Value = 255
$ws4.Range("F23").Value = 1309
$ws4.Range("F24").Value = 1281
$ws4.Range("F26").Value = 1300
$ws4.Range("F35").Value = 151
$ws4.Range("F38").Value = 585
$ws4.Range("F40").Value = 126
$ws4.Range("F41").Value = 133
$ws4.Range("F42").Value = 450
$ws4.Range("F43").Value = 636
$ws4.Range("F45").Value = 178
$ws4.Range("F46").Value = 87
$ws4.Range("F47").Value = 35
$ws4.Range("F49").Value = 13
